$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header cells (e.g. H1: bold font, thin border, centered/top aligned).
$ws.Range("I1").Value = "I0"
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats

# Add the corresponding data values in row 2 (plain, unstyled numeric cells).
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
